$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# D7: rule action text changes from setPrecioOro to setValorAvaluo
$ws.Range("D7").Value = "garantia.setValorAvaluo(Double.valueOf(`$param));"

# Numeric value updates in column D
$ws.Range("D12").Value = 18.5
$ws.Range("D13").Value = 18
$ws.Range("D15").Value = 17

# Update selection to D16 (matches the sheetView selection change in the diff)
$ws.Range("D16").Select()
